{"js": "// Applies the Project Write-Up Chatapredu edits:\n// 1. \"The project may include\" -> \"This project may include\"\n// 2. Tech stack sentence expanded with JavaScript/TypeScript, Redux, Cognito,\n//    Serverless Framework, AWS Amplify; \"and\" dropped before \"AWS API Gateway\"\n// 3. \"Frontend\" -> \"Front-end\"\n// 4. \"Backend\" -> \"Back-end\"\n// 5. \"storing user tables, storing group tables\" -> \"storing all persistent\n//    data using a comprehensive table\"\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"The project may include additional features\",\n  \"This project may include additional features\"\n);\n\nawait replaceOnce(\n  \"React Native, AWS Lambda, PostgreSQL, Node.js, and AWS API Gateway, and AWS S3.\",\n  \"JavaScript/TypeScript, React Native supported with Redux, AWS Lambda, PostgreSQL, Node.js, AWS API Gateway, Cognito, Serverless Framework, AWS Amplify, and AWS S3.\"\n);\n\nawait replaceOnce(\n  \"Frontend responsibilities include:\",\n  \"Front-end responsibilities include:\"\n);\n\nawait replaceOnce(\n  \"Backend responsibilities include:\",\n  \"Back-end responsibilities include:\"\n);\n\nawait replaceOnce(\n  \"Database responsibilities include: storing user tables, storing group tables, communication with the backend,\",\n  \"Database responsibilities include: storing all persistent data using a comprehensive table, communication with the backend,\"\n);\n", "ps1": "# Applies the Project Write-Up Chatapredu edits:\n# 1. \"The project may include\" -> \"This project may include\"\n# 2. Tech stack sentence expanded with JavaScript/TypeScript, Redux, Cognito,\n#    Serverless Framework, AWS Amplify; \"and\" dropped before \"AWS API Gateway\"\n# 3. \"Frontend\" -> \"Front-end\"\n# 4. \"Backend\" -> \"Back-end\"\n# 5. \"storing user tables, storing group tables\" -> \"storing all persistent\n#    data using a comprehensive table\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\nReplace-Once \"The project may include additional features\" \"This project may include additional features\"\n\nReplace-Once \"React Native, AWS Lambda, PostgreSQL, Node.js, and AWS API Gateway, and AWS S3.\" \"JavaScript/TypeScript, React Native supported with Redux, AWS Lambda, PostgreSQL, Node.js, AWS API Gateway, Cognito, Serverless Framework, AWS Amplify, and AWS S3.\"\n\n# \"Frontend\"/\"Backend\" each sit in their own run immediately followed by a\n# proofErr-wrapped \"include:\" run; only touch the word itself so the\n# gramStart/gramEnd pair around \"include:\" stays balanced.\nReplace-Once \"Frontend responsibilities \" \"Front-end responsibilities \"\n\nReplace-Once \"Backend responsibilities \" \"Back-end responsibilities \"\n\nReplace-Once \"Database responsibilities include: storing user tables, storing group tables, communication with the backend,\" \"Database responsibilities include: storing all persistent data using a comprehensive table, communication with the backend,\"\n"}
